$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Likes / Dislikes stay the same text, but the "low reward memes"
# helper column (G) and its companion emoji column (H) go away.
$ws.Range("B1").Value = "Likes"
$ws.Range("C1").Value = "Dislikes"
$ws.Range("G1:H1").ClearContents()

# Column A used to hold "memeNNN" labels; the feedback page now wants the
# plain numeric meme id (101-208) instead.
for ($row = 2; $row -le 109; $row++) {
    $ws.Cells.Item($row, 1).Value = 100 + ($row - 1)
}

# Match the recorded selection left behind by the edit.
[void]$ws.Range("G1").Select()
